# Generate Report for Handoff
#
# The four "Ready for handoff" rows (4-7: 4ac08ebc..., a38c41b4..., d88b5331...,
# ed969ac8...) on each language sheet just had a fresh handoff xliff generated:
#   - Priority flips from "low" to "ht" (now hand-translated/in-flight)
#   - Latest Handoff Datetime is bumped to the new generation timestamp
# The Overview sheet's "Latest HO Xliff Generate Date" column mirrors the
# newest per-language handoff time (de-de's, since it is later than zh-cn's),
# so it is refreshed to match as well.

$wb = $excel.ActiveWorkbook

$zhRows = @(4, 5, 6, 7)
$deRows = @(4, 5, 6, 7)

$zh = $wb.Worksheets.Item("zh-cn")
foreach ($r in $zhRows) {
    $zh.Cells.Item($r, 5).Value = "ht"
    $zh.Cells.Item($r, 8).Value = "2016-08-23 12:30:39"
}

$de = $wb.Worksheets.Item("de-de")
foreach ($r in $deRows) {
    $de.Cells.Item($r, 5).Value = "ht"
    $de.Cells.Item($r, 8).Value = "2016-08-23 12:30:43"
}

$ov = $wb.Worksheets.Item("Overview")
foreach ($r in @(4, 5, 6, 7)) {
    $ov.Cells.Item($r, 7).Value = "2016-08-23 12:30:43"
}
